# Update column F (dSF) values for specific rows per repull/push of data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F7").Value = -3
$ws.Range("F10").Value = 5
$ws.Range("F16").Value = -3
$ws.Range("F18").Value = -11
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = -2
$ws.Range("F22").Value = -6
$ws.Range("F23").Value = 4
